# Update both the "crosstab" (numeric) and "annot" (text/inline-string)
# worksheets with the corrected values for the 04dec2025 column (E), and
# also fix the 03dec2025 column (D) for one row.
#
# Column E corresponds to header "04dec2025" and column D to "03dec2025".
#
# Changes (row -> new value):
#   Row 3  (CAMACHO LINARES JUDITH ARLETT):        E  0  -> 8
#   Row 4  (CONTRERAS VALDERRAMA JULIA ALEJANDRA):  E  0  -> 10
#   Row 5  (FERNANDEZ VALDERAS ERNESTO ALI):        E  8  -> 21
#   Row 6  (GUTIERREZ CARLOS TERESA DE JESUS):      E  0  -> 7
#   Row 8  (HUMPIRE CASTILLO IRWIN DEIMER):         D  7  -> 6 , E 0 -> 3
#   Row 12 (SEVERINO AVALOS MARJORIE ISABEL):       E  0  -> 8
#   Row 13 (VALLE MAGALLAN EDUAR):                  E  0  -> 12
#   Row 14 (ZAVALETA MANAY JORGE LUIS):             E  0  -> 3
#   Row 15 (ZEVALLOS PACHECO ZOILA XIMENA):         E  0  -> 7
#
# The "crosstab" sheet keeps these as numbers; the "annot" sheet keeps the
# same figures but stored as text (as the rest of that sheet already is),
# so the cell's number format is forced to Text before the value is written
# -- otherwise a numeric-looking string like "8" would be auto-converted
# back into a number by Excel.

$wb = $excel.ActiveWorkbook

$crosstab = $wb.Worksheets.Item("crosstab")
$annot = $wb.Worksheets.Item("annot")

# Map of row -> new column E value (numbers)
$eUpdates = @{
    3 = 8
    4 = 10
    5 = 21
    6 = 7
    8 = 3
    12 = 8
    13 = 12
    14 = 3
    15 = 7
}

foreach ($row in $eUpdates.Keys) {
    $value = $eUpdates[$row]

    # Numeric crosstab sheet
    $crosstab.Range("E$row").Value = $value

    # Text/annotation sheet stores the same number as text
    $annot.Range("E$row").NumberFormat = "@"
    $annot.Range("E$row").Value = "$value"
}

# Row 8 also has a change in column D (7 -> 6)
$crosstab.Range("D8").Value = 6
$annot.Range("D8").NumberFormat = "@"
$annot.Range("D8").Value = "6"
